# Apply updated absenteeism records to rows 2-11 of the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: RowNumber = @(A_id, B_name, C_department, D_reason, E_hours, F_date_serial, G_salary)
$data = @{
    2  = @(28320, "Rodrigo Nogueira",      "Atendimento ao Cliente", "Outros",              2, 45101, 12055.55)
    3  = @(90269, "Dr. Thales Aragão",     "Marketing",              "Viagem de negócios",  5, 45101, 3100.08)
    4  = @(51475, "Ryan Ferreira",         "Vendas",                 "Problemas pessoais",  5, 45104, 7401.55)
    5  = @(2527,  "Melissa da Paz",        "Vendas",                 "Outros",              3, 45102, 10589.57)
    6  = @(14799, "Miguel das Neves",      "Marketing",              "Outros",              1, 45106, 6754.14)
    7  = @(60452, "Sabrina Alves",         "Operações",              "Outros",              2, 45080, 11793.24)
    8  = @(32595, "Pietra Santos",         "Operações",              "Doença",              4, 45085, 6612.3)
    9  = @(96191, "Bernardo Fernandes",    "Atendimento ao Cliente", "Outros",              4, 45104, 6490.82)
    10 = @(99016, "Sofia Castro",          "Vendas",                 "Consulta médica",     1, 45095, 9064.68)
    11 = @(27612, "Luna Melo",             "Atendimento ao Cliente", "Outros",              4, 45099, 3262.61)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]

    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
    $ws.Cells.Item($row, 6).Value = $values[5]
    $ws.Cells.Item($row, 7).Value = $values[6]
}
